# Fruta / hortaliza, semanal
# Insert a new weekly data point for "Zanahoria" (Terminal La Palmera de La Serena)
# as a new row 442, pushing the existing rows 442:461 down to 443:462.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 442 (shifts 442:461 -> 443:462)
$ws.Rows.Item(442).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A442").Value = 8
$ws.Range("B442").Value = "Terminal La Palmera de La Serena"
$ws.Range("C442").Value = "Coquimbo"
$ws.Range("D442").Value = 44939
$ws.Range("E442").Value = 4
$ws.Range("F442").Value = 100114013
$ws.Range("G442").Value = "Zanahoria"
$ws.Range("H442").Value = "Sin especificar"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 400
$ws.Range("K442").Value = 6000
$ws.Range("L442").Value = 7000
$ws.Range("M442").Value = 6500
$ws.Range("N442").Value = "$/saco 20 kilos"
$ws.Range("O442").Value = "Provincia del Elquí"
$ws.Range("P442").Value = 325
$ws.Range("Q442").Value = 20
$ws.Range("R442").Value = "Hortaliza"

# Keep the date column formatted consistently with the rest of the column
$ws.Range("D442").NumberFormat = "YYYY-MM-DD HH:MM:SS"
